# EliteSupportDatasheet.pptx - proofreading fix on the Elite Support
# hours-of-operation table (slide 4): add missing end-of-sentence
# periods to the two footnote lines in the merged "Americas" cell of
# "Table 6" (shape/graphicFrame id 25).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Locate the table shape by its (stable) shape Id rather than a
# positional index - the slide has two shapes named "Table 6"
# (ids 25 and 111); the footnote lives in the one with id 25.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Id -eq 25) {
        $shp = $cand
        break
    }
}

$tbl = $shp.Table

# The footnote text lives in the first (left-most) cell of row 3,
# which spans all 4 columns. It holds three paragraphs:
#   1) "Language support is only available in English and Japanese "
#   2) (blank)
#   3) " 1 P2, P3, P4 cases are limited to business hours only in Japan"
$cell = $tbl.Cell(3, 1)
$tr = $cell.Shape.TextFrame.TextRange

$nbsp = [char]0x00A0

$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "Language support is only available in English and Japanese." + $nbsp

$para3 = $tr.Paragraphs(3, 1)
$para3.Text = $nbsp + "1 P2, P3, P4 cases are limited to business hours only in Japan."
